# Inspire address harmonised postal code 08011:
# add a new SOURCE/TAXONOMY mapping row ("-" -> "None") to the
# conditionOfConstruction sheet.

$wb = $excel.ActiveWorkbook

$wsCondition = $wb.Worksheets.Item("conditionOfConstruction")
$wsMainUse   = $wb.Worksheets.Item("mainUse")

$wsCondition.Range("A8").Value = "-"
$wsCondition.Range("B8").Value = "None"

# Leave the selection on mainUse where it last was (matches author's
# recorded cursor position) before finally landing on conditionOfConstruction.
$wsMainUse.Select() | Out-Null
$wsMainUse.Range("A17").Select() | Out-Null

$wsCondition.Select() | Out-Null
$wsCondition.Range("C13").Select() | Out-Null
